# Change range of M to also use survey data
#
# 1) descriptions sheet: E10..E19 ("M50".."M150") descriptions become "external"
# 2) descriptions sheet: 10 new rows (E20..E29 -> "internal") inserted before the
#    F-row block, pushing F0/F1/F2/L10/L30/L31 down by 10 rows
# 3) scenarios sheet: rows 40-49 and 55-64 switch from case A10/L10 ("fish age" /
#    "fish length") to A30/L30 ("fish & surv length") for columns A and D

$wb = $excel.ActiveWorkbook
$wsScenarios = $wb.Worksheets.Item("scenarios")
$wsDescriptions = $wb.Worksheets.Item("descriptions")

# --- 1) descriptions!D16:D25 (cases E10..E19) -> "external" ---
for ($r = 16; $r -le 25; $r++) {
    $wsDescriptions.Cells.Item($r, 4).Value = "external"
}

# --- 2) insert 10 rows before row 26, fill with E20..E29 / internal ---
$wsDescriptions.Rows("26:35").Insert()

for ($i = 0; $i -le 9; $i++) {
    $r = 26 + $i
    $num = 20 + $i
    $wsDescriptions.Cells.Item($r, 1).Value = "E"
    $wsDescriptions.Cells.Item($r, 2).Value = $num
    $wsDescriptions.Cells.Item($r, 3).Formula = "=CONCATENATE(A" + $r + ",B" + $r + ")"
    $wsDescriptions.Cells.Item($r, 4).Value = "internal"
}

# Selection ends up parked on descriptions!E25 (matching the saved view state),
# while the scenarios sheet stays the active/selected tab.
$wsDescriptions.Activate()
$wsDescriptions.Range("E25").Select()
$wsScenarios.Activate()

# --- 3) scenarios!A/D for rows 40-49 and 55-64: A10/L10 -> A30/L30 ---
$targetRows = @(40,41,42,43,44,45,46,47,48,49,55,56,57,58,59,60,61,62,63,64)
foreach ($r in $targetRows) {
    $wsScenarios.Cells.Item($r, 1).Value = "A30"
    $wsScenarios.Cells.Item($r, 4).Value = "L30"
}
